# Apply the "subtest added (alphabet writing fluency); fix allows commas in
# student specific info; fix to example template" edit described in the
# commit diff:
#   - every "Student Specific Information for ... [subtest]" sentence in
#     column C (rows 5-41) gains a trailing period (Google-Sheets-safe CSV
#     fix so a stray comma inside the text no longer looks "unfinished"),
#   - row 33 (Alphabet Writing Fluency) now documents that the subtest was
#     NOT administered, instead of the generic "Student Specific
#     Information..." placeholder,
#   - the big "Conclusions:" example text (C43) is reworded slightly,
#   - the view scrolls down to/selects C33, the row that changed meaning.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(5, 3).Value  = "Student Specific Information for the Oral Language Composite."
$ws.Cells.Item(6, 3).Value  = "Student Specific Information for the Listening Comprehension subtest."
$ws.Cells.Item(7, 3).Value  = "Student Specific Information for the Receptive Vocabulary subtest."
$ws.Cells.Item(8, 3).Value  = "Student Specific Information for the Oral Discourse Comprehension."
$ws.Cells.Item(9, 3).Value  = "Student Specific Information for the Oral Expression subtest."
$ws.Cells.Item(10, 3).Value = "Student Specific Information for the Expressive Vocabulary subtest."
$ws.Cells.Item(11, 3).Value = "Student Specific Information for the Oral Word Fluency subtest."
$ws.Cells.Item(12, 3).Value = "Student Specific Information for the Sentence Repetition subtest."
$ws.Cells.Item(13, 3).Value = "Student Specific Information for the Phonological Processing Composite."
$ws.Cells.Item(14, 3).Value = "Student Specific Information for the Pseudoword Decoding subtest."
$ws.Cells.Item(15, 3).Value = "Student Specific Information for Phonemic Proficiency subtest."
$ws.Cells.Item(16, 3).Value = "Student Specific Information for the Orthographic Processing Composite."
$ws.Cells.Item(17, 3).Value = "Student Specific Information for the Orthographic Fluency subtest."
$ws.Cells.Item(18, 3).Value = "Student Specific Information for the Spelling subtest."
$ws.Cells.Item(19, 3).Value = "Student Specific Information for the Orthographic Choice subtest."
$ws.Cells.Item(20, 3).Value = "Student Specific Information for the Reading Composite."
$ws.Cells.Item(21, 3).Value = "Student Specific Information for the Word Reading subtest."
$ws.Cells.Item(22, 3).Value = "Student Specific Information for the Reading Comprehension subtest."
$ws.Cells.Item(23, 3).Value = "Student Specific Information for the Basic Reading (and Decoding) Composite."
$ws.Cells.Item(24, 3).Value = "Student Specific Information for the Reading Fluency Composite."
$ws.Cells.Item(25, 3).Value = "Student Specific Information for the Oral Reading Fluency subtest."
$ws.Cells.Item(26, 3).Value = "Student Specific Information for the Decoding Fluency subtest."
$ws.Cells.Item(27, 3).Value = "Student Specific Information for the Written Expression Composite."
$ws.Cells.Item(28, 3).Value = "Student Specific Information for the Sentence Composition subtest."
$ws.Cells.Item(29, 3).Value = "Student Specific Information for the Essay Composition subtest."
$ws.Cells.Item(30, 3).Value = "Student Specific Information for the Sentence Building subtest."
$ws.Cells.Item(31, 3).Value = "Student Specific Information for the Sentence Combining subtest."
$ws.Cells.Item(32, 3).Value = "Student Specific Information for the Writing Fluency composite."

# Newly-documented "not administered" subtest (Alphabet Writing Fluency).
$ws.Cells.Item(33, 3).Value = "Student Name did not have the Alphabet Writing Fluency subtest administered because, [insert description here]."

$ws.Cells.Item(34, 3).Value = "Student Specific Information for Sentence Writing Fluency subtest."
$ws.Cells.Item(35, 3).Value = "Student Specific Information for the Mathematics Composite."
$ws.Cells.Item(36, 3).Value = "Student Specific Information for the Math Problem Solving subtest."
$ws.Cells.Item(37, 3).Value = "Student Specific Information for the Numerical Operations subtest."
$ws.Cells.Item(38, 3).Value = "Student Specific Information for the Math Fluency Composite."
$ws.Cells.Item(39, 3).Value = "Student Specific Information for the Math Fluency–Addition subtest."
$ws.Cells.Item(40, 3).Value = "Student Specific Information for the Math Fluency–Subtraction subtest."
$ws.Cells.Item(41, 3).Value = "Student Specific Information for Math Fluency–Multiplication subtest."

# Reworded "Conclusions:" example/instructions cell.
$conclusion  = "This is where we can place all of our conclusions about Student-Name's performance on the exam. This way, this single CSV file will contain all of the important results of the student at-a-glance, which I have been told will be a very useful document for the Special Education Teacher as before a meeting, the teacher can look at all of the information in one place (especially when they have so many students on their caseload), which is excellent, as any functionality which will improve the work-flow of a teacher, or simply improve the quality of work-life of a teacher, is certainly worth the (minor) efforts to include this functionality in the report generator.  With Word Wrapping on a cell, the conclusion can be very readable indeed within a spreadsheet application. With this addition of adding concluding remarks, all of the student specific data and all of the teacher's observations/recordings are in one single document; quite glorious indeed!"
$conclusion += "`n"
$conclusion += "We can even add a bulleted list into this conclusion (a copy/paste of a bullet symbol is required). Example: The following list presents my recomendations based on this test for Student-Name:`n"
$conclusion += "   " + [char]0x2022 + " Recomendation list item #1`n"
$conclusion += "   " + [char]0x2022 + " Recomendation list item #2`n"
$conclusion += "   " + [char]0x2022 + " Recomendation list item #3"

$ws.Cells.Item(43, 3).Value = $conclusion

# Move the visible window so it is scrolled to/selecting the row that
# changed meaning (Alphabet Writing Fluency, C33) instead of the old C25.
$ws.Range("C33").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
